$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93 (shifts old rows 93:113 down to 94:114)
$ws.Rows.Item(93).Insert()

# New row 93: backfilled date 2021-02-13 (44235), with 0 new cases that day.
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 0
$ws.Range("C93").Value = 2
$ws.Range("D93").Value = 53.53319057815846

# Row 92 (44234) rolling 7-day sum is recomputed now that 44235 sits between
# 44234 and 44236.
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 53.53319057815846

# Row 112 (44254) now has an updated rolling 7-day sum too.
$ws.Range("C112").Value = 15
$ws.Range("D112").Value = 401.4989293361884

# Append a new row 115 for the following day 2021-03-18 (44257), 3 new cases,
# rolling sums not yet computed (left blank) just like the other recent days.
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 3
# Match the date display format used by the rest of column A.
$ws.Range("A115").NumberFormat = "YYYY-MM-DD HH:MM:SS"
